$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: phone "09876543" (text, leading zero) -> numeric 9876543
$ws.Range("A25").Value = 9876543

# New row 26: re-add the customer record with the leading-zero phone
# "09876543" as text (so the leading zero is preserved), points reset to 0.
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "09876543"
$ws.Range("A26").Style = "Normal"

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Style = "Normal"

$ws.Range("C26").Value = 0
